# Apply data refresh to "展览" and "全部类型" sheets: update 想去人数 (F) / 最低票价 (G)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2202
$ws1.Range("G2").Value = 60
$ws1.Range("F3").Value = 629
$ws1.Range("F4").Value = 1588
$ws1.Range("F5").Value = 7395
$ws1.Range("F7").Value = 187

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2202
$ws4.Range("F3").Value = 629
$ws4.Range("F4").Value = 1588
$ws4.Range("F5").Value = 7395
$ws4.Range("F7").Value = 187
